$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 99.28451533333333
$ws.Range("H2").Value = 297.853546
$ws.Range("I2").Value = 0.02270354261926982
$ws.Range("J2").Value = 0.02270354261926982
$ws.Range("M2").Value = 7.407905
$ws.Range("N2").Value = 22.223715
$ws.Range("O2").Value = 0.1577242380174723
$ws.Range("P2").Value = 0.1577242380174723
$ws.Range("Q2").Value = 735.4902575603766
$ws.Range("R2").Value = 6619.41231804339
$ws.Range("S2").Value = 0.00358089895992154
$ws.Range("T2").Value = 0.00358089895992154

$ws.Range("G3").Value = 99.28451533333333
$ws.Range("H3").Value = 297.853546
$ws.Range("I3").Value = 0.02270354261926982
$ws.Range("J3").Value = 0.02270354261926982
$ws.Range("O3").Value = 0.3510414535684271
$ws.Range("P3").Value = 0.3510414535684271
$ws.Range("Q3").Value = 1636.955564627994
$ws.Range("R3").Value = 14732.60008165194
$ws.Range("S3").Value = 0.007969884602221214
$ws.Range("T3").Value = 0.007969884602221214

$ws.Range("G4").Value = 99.28451533333333
$ws.Range("H4").Value = 297.853546
$ws.Range("I4").Value = 0.02270354261926982
$ws.Range("J4").Value = 0.02270354261926982
$ws.Range("M4").Value = 5.464566666666666
$ws.Range("N4").Value = 16.3937
$ws.Range("O4").Value = 0.1163479571613943
$ws.Range("P4").Value = 0.1163479571613943
$ws.Range("Q4").Value = 542.5468530066888
$ws.Range("R4").Value = 4882.9216770602
$ws.Range("S4").Value = 0.002641510804078694
$ws.Range("T4").Value = 0.002641510804078694

$ws.Range("G5").Value = 99.28451533333333
$ws.Range("H5").Value = 297.853546
$ws.Range("I5").Value = 0.02270354261926982
$ws.Range("J5").Value = 0.02270354261926982
$ws.Range("M5").Value = 17.60745533333333
$ws.Range("N5").Value = 52.822366
$ws.Range("O5").Value = 0.3748863512527063
$ws.Range("P5").Value = 0.3748863512527063
$ws.Range("Q5").Value = 1748.147669023315
$ws.Range("R5").Value = 15733.32902120984
$ws.Range("S5").Value = 0.008511248253048373
$ws.Range("T5").Value = 0.008511248253048373

$ws.Range("I6").Value = 0.9171714767027319
$ws.Range("J6").Value = 0.9171714767027318
$ws.Range("M6").Value = 7.407905
$ws.Range("N6").Value = 22.223715
$ws.Range("O6").Value = 0.1577242380174723
$ws.Range("P6").Value = 0.1577242380174723
$ws.Range("Q6").Value = 29712.13334145375
$ws.Range("R6").Value = 267409.2000730838
$ws.Range("S6").Value = 0.1446601722942982
$ws.Range("T6").Value = 0.1446601722942982

$ws.Range("I7").Value = 0.9171714767027319
$ws.Range("J7").Value = 0.9171714767027318
$ws.Range("O7").Value = 0.3510414535684271
$ws.Range("P7").Value = 0.3510414535684271
$ws.Range("S7").Value = 0.3219652083532278
$ws.Range("T7").Value = 0.3219652083532278

$ws.Range("I8").Value = 0.9171714767027319
$ws.Range("J8").Value = 0.9171714767027318
$ws.Range("M8").Value = 5.464566666666666
$ws.Range("N8").Value = 16.3937
$ws.Range("O8").Value = 0.1163479571613943
$ws.Range("P8").Value = 0.1163479571613943
$ws.Range("Q8").Value = 21917.65869746757
$ws.Range("R8").Value = 197258.9282772081
$ws.Range("S8").Value = 0.1067110276810622
$ws.Range("T8").Value = 0.1067110276810622

$ws.Range("I9").Value = 0.9171714767027319
$ws.Range("J9").Value = 0.9171714767027318
$ws.Range("M9").Value = 17.60745533333333
$ws.Range("N9").Value = 52.822366
$ws.Range("O9").Value = 0.3748863512527063
$ws.Range("P9").Value = 0.3748863512527063
$ws.Range("Q9").Value = 70621.18921175301
$ws.Range("R9").Value = 635590.702905777
$ws.Range("S9").Value = 0.3438350683741436
$ws.Range("T9").Value = 0.3438350683741436

$ws.Range("G10").Value = 249.2612966666667
$ws.Range("H10").Value = 747.7838899999999
$ws.Range("I10").Value = 0.05699896356653876
$ws.Range("J10").Value = 0.05699896356653875
$ws.Range("M10").Value = 7.407905
$ws.Range("N10").Value = 22.223715
$ws.Range("O10").Value = 0.1577242380174723
$ws.Range("P10").Value = 0.1577242380174723
$ws.Range("Q10").Value = 1846.504005883483
$ws.Range("R10").Value = 16618.53605295135
$ws.Range("S10").Value = 0.00899011809631799
$ws.Range("T10").Value = 0.008990118096317992

$ws.Range("G11").Value = 249.2612966666667
$ws.Range("H11").Value = 747.7838899999999
$ws.Range("I11").Value = 0.05699896356653876
$ws.Range("J11").Value = 0.05699896356653875
$ws.Range("O11").Value = 0.3510414535684271
$ws.Range("P11").Value = 0.3510414535684271
$ws.Range("Q11").Value = 4109.700946365995
$ws.Range("R11").Value = 36987.30851729395
$ws.Range("S11").Value = 0.02000899902229158
$ws.Range("T11").Value = 0.02000899902229158

$ws.Range("G12").Value = 249.2612966666667
$ws.Range("H12").Value = 747.7838899999999
$ws.Range("I12").Value = 0.05699896356653876
$ws.Range("J12").Value = 0.05699896356653875
$ws.Range("M12").Value = 5.464566666666666
$ws.Range("N12").Value = 16.3937
$ws.Range("O12").Value = 0.1163479571613943
$ws.Range("P12").Value = 0.1163479571613943
$ws.Range("Q12").Value = 1362.104973054778
$ws.Range("R12").Value = 12258.944757493
$ws.Range("S12").Value = 0.006631712971283525
$ws.Range("T12").Value = 0.006631712971283525

$ws.Range("G13").Value = 249.2612966666667
$ws.Range("H13").Value = 747.7838899999999
$ws.Range("I13").Value = 0.05699896356653876
$ws.Range("J13").Value = 0.05699896356653875
$ws.Range("M13").Value = 17.60745533333333
$ws.Range("N13").Value = 52.822366
$ws.Range("O13").Value = 0.3748863512527063
$ws.Range("P13").Value = 0.3748863512527063
$ws.Range("Q13").Value = 4388.857147387082
$ws.Range("R13").Value = 39499.71432648374
$ws.Range("S13").Value = 0.02136813347664565
$ws.Range("T13").Value = 0.02136813347664566

$ws.Range("G14").Value = 13.67033766666667
$ws.Range("H14").Value = 41.011013
$ws.Range("I14").Value = 0.003126017111459632
$ws.Range("J14").Value = 0.003126017111459632
$ws.Range("M14").Value = 7.407905
$ws.Range("N14").Value = 22.223715
$ws.Range("O14").Value = 0.1577242380174723
$ws.Range("P14").Value = 0.1577242380174723
$ws.Range("Q14").Value = 101.2685627525883
$ws.Range("R14").Value = 911.4170647732949
$ws.Range("S14").Value = 0.0004930486669345503
$ws.Range("T14").Value = 0.0004930486669345503

$ws.Range("G15").Value = 13.67033766666667
$ws.Range("H15").Value = 41.011013
$ws.Range("I15").Value = 0.003126017111459632
$ws.Range("J15").Value = 0.003126017111459632
$ws.Range("O15").Value = 0.3510414535684271
$ws.Range("P15").Value = 0.3510414535684271
$ws.Range("Q15").Value = 225.3899839130369
$ws.Range("R15").Value = 2028.509855217332
$ws.Range("S15").Value = 0.001097361590686565
$ws.Range("T15").Value = 0.001097361590686565

$ws.Range("G16").Value = 13.67033766666667
$ws.Range("H16").Value = 41.011013
$ws.Range("I16").Value = 0.003126017111459632
$ws.Range("J16").Value = 0.003126017111459632
$ws.Range("M16").Value = 5.464566666666666
$ws.Range("N16").Value = 16.3937
$ws.Range("O16").Value = 0.1163479571613943
$ws.Range("P16").Value = 0.1163479571613943
$ws.Range("Q16").Value = 74.70247153534444
$ws.Range("R16").Value = 672.3222438180999
$ws.Range("S16").Value = 0.0003637057049698908
$ws.Range("T16").Value = 0.0003637057049698908

$ws.Range("G17").Value = 13.67033766666667
$ws.Range("H17").Value = 41.011013
$ws.Range("I17").Value = 0.003126017111459632
$ws.Range("J17").Value = 0.003126017111459632
$ws.Range("M17").Value = 17.60745533333333
$ws.Range("N17").Value = 52.822366
$ws.Range("O17").Value = 0.3748863512527063
$ws.Range("P17").Value = 0.3748863512527063
$ws.Range("Q17").Value = 240.6998598574176
$ws.Range("R17").Value = 2166.298738716758
$ws.Range("S17").Value = 0.001171901148868626
$ws.Range("T17").Value = 0.001171901148868626
